$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.524.19"
$ws.Range("E2").Value = "  -4.53%  "
$ws.Range("D3").Value = "3.285.41"
$ws.Range("E3").Value = "  -6.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.281.97"
$ws.Range("E8").Value = "  -6.99%  "
$ws.Range("E9").Value = "  -10.54%  "
$ws.Range("E10").Value = "  -13.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -16.17%  "
$ws.Range("E14").Value = "  -10.82%  "
$ws.Range("D15").Value = "3.819.20"
$ws.Range("E15").Value = "  -6.88%  "
$ws.Range("D16").Value = "67.626.88"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "3.285.97"
$ws.Range("E17").Value = "  -7.08%  "
$ws.Range("E18").Value = "  -13.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "537.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -11.45%  "
$ws.Range("E20").Value = "  -6.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.768"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -11.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.41%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.52%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.13%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -15.28%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.15%  "
$ws.Range("E32").Value = "  -11.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "546.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -17.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.81"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -14.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0465"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.18%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "53.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0869"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -15.50%  "
$ws.Range("E41").Value = "  -9.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "2.952.64"
$ws.Range("E43").Value = "  -11.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.271"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -12.64%  "
$ws.Range("D45").Value = "0.0₃0598"
$ws.Range("E45").Value = "  -17.56%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.17%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -15.70%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -18.31%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("E51").Value = "  -11.96%  "
